$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that were marked as "Killed" / "survived" (or "Survived") after running
# the Botium test cases against each mutant. For each of these rows we mark
# column B as "Done" and column E with the mutation-testing outcome.
$results = @{
    4  = "survived"
    8  = "Killed"
    13 = "survived"
    19 = "survived"
    20 = "survived"
    21 = "Killed"
    23 = "Killed"
    24 = "survived"
    26 = "survived"
    27 = "survived"
    31 = "Killed"
    35 = "Killed"
    43 = "survived"
    80 = "Survived"
}

$rowsInOrder = 4, 8, 13, 19, 20, 21, 23, 24, 26, 27, 31, 35, 43, 80

foreach ($r in $rowsInOrder) {
    $status = $results[$r]
    $ws.Range("B" + $r).Value = "Done"
    $ws.Range("E" + $r).Value = $status
}

# Manually entered Botium test-case count.
$ws.Range("I3").Value = 93

# Normalize the formatting of every non-highlighted data row (A2:F94) so it
# matches the plain bordered style already used elsewhere on the sheet
# (instead of the redundant "border + explicit no-fill" style).
$ws.Range("A95").Copy()
for ($r = 2; $r -le 94; $r++) {
    $cell = $ws.Range("A" + $r)
    $color = $cell.Interior.Color()
    if ($color -eq 16777215) {
        $rng = $ws.Range("A" + $r + ":F" + $r)
        $rng.PasteSpecial(-4122)
    }
}

$ws.Range("B81").Select()
